$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$m.HeadersFooters.DateAndTime.Text = "08-09-2019"
Write-Host "Text=$($m.HeadersFooters.DateAndTime.Text)"
